$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Report" to "Sheet1"
$ws.Name = "Sheet1"

# Update membership count column (B) with corrected numeric values
$ws.Range("B2").Value = 1100
$ws.Range("B3").Value = 1000
$ws.Range("B4").Value = 125
$ws.Range("B5").Value = 300
$ws.Range("B6").Value = 554
